# Finalise Technical Manual and Agent Manual
# Applies the "Use Cases" workbook update:
#   - Event Table (sheet1): inserts a new "Assign user" event row (new row 5)
#     and a new "Agent has no users to assign to themselves" event row (new
#     row 9), renumbers the activity steps in the existing rows, and tweaks
#     a few trigger / response strings.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row formatting helper values (match existing table look: bordered,
#     centered, vertically centered, wrapped text) ---

function Format-EventRow($range) {
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4108
    $range.WrapText = $true
}

# 1. Insert the new "Assign user" event as row 5 (pushes the existing
#    "User details received" / "Query received" / "Answer to query found"
#    rows down by one).
$ws1.Rows.Item(5).Insert()
$ws1.Rows.Item(5).RowHeight = 50.1

$ws1.Range("A5").Value = "Agent sees unassigned user in open queue"
$ws1.Range("B5").Value = "State"
$ws1.Range("C5").Value = "(Triggered from Use Case 1 or 2)"
$ws1.Range("E5").Value = "3. Assign user"
$ws1.Range("F5").Value = "UserID, AssignedTo, Status"
Format-EventRow $ws1.Range("A5:G5")

# 2. Update the (now shifted) "User details received from user" row (row 6)
$ws1.Range("C6").Value = "(Triggered from Use Case 3)"
$ws1.Range("E6").Value = "4. Update user details"
$ws1.Range("F6").Value = "UserID, Query content"

# 3. Update the (now shifted) "Query received from user" row (row 7)
$ws1.Range("C7").Value = "(Triggered from Use Case 3)"
$ws1.Range("E7").Value = "5. Look up answer"
$ws1.Range("F7").Value = "Query content, Time, Location, Audio, Image"

# 4. Update the (now shifted) "Answer to query found" row (row 8)
$ws1.Range("C8").Value = "(Triggered from Use Case 5)"
$ws1.Range("E8").Value = "6. Send response"
$ws1.Range("F8").Value = "Answer"
Format-EventRow $ws1.Range("A8:G8")

# 5. Insert the new "Agent has no users to assign to themselves" event as
#    row 9 (after the "Answer to query found" row, before the blank rows).
$ws1.Rows.Item(9).Insert()
$ws1.Rows.Item(9).RowHeight = 50.1

$ws1.Range("A9").Value = "Agent has no users to assign to themselves"
$ws1.Range("B9").Value = "State"
$ws1.Range("C9").Value = "(Triggered from Use Case 3)"
$ws1.Range("E9").Value = "7. Request user release"
$ws1.Range("F9").Value = "AgentID"
Format-EventRow $ws1.Range("A9:G9")
$ws1.Range("C9").Borders.LineStyle = 1
$ws1.Range("C9").HorizontalAlignment = -4108
$ws1.Range("C9").VerticalAlignment = -4108
$ws1.Range("C9").WrapText = $true

# --- View state: scroll back to the top and select the finished table ---
$ws1.Activate()
$win = $ws1.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws1.Range("A1:G9").Select()

Write-Host "Event Table updated"
